# Add "2022-Q3" worksheet with fund holding data, inserted right after
# "总计" (and before "2022-Q2"), and record the new quarter in the "总计"
# summary sheet (new row 2, with existing rows shifting down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a row for 2022-Q3 data,
#    shifting the 2022-Q2 / 2021-Q3 / 2021-Q2 rows down by one.
#    (Column A is a plain running index 0,1,2,3 and is left alone except
#    for the brand new row 5.)
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$b4 = $summary.Range("B4").Value()
$c4 = $summary.Range("C4").Value()
$d4 = $summary.Range("D4").Value()
$b3 = $summary.Range("B3").Value()
$c3 = $summary.Range("C3").Value()
$d3 = $summary.Range("D3").Value()
$b2 = $summary.Range("B2").Value()
$c2 = $summary.Range("C2").Value()
$d2 = $summary.Range("D2").Value()

$summary.Range("B5").Value = $b4
$summary.Range("C5").Value = $c4
$summary.Range("D5").Value = $d4

$summary.Range("B4").Value = $b3
$summary.Range("C4").Value = $c3
$summary.Range("D4").Value = $d3

$summary.Range("B3").Value = $b2
$summary.Range("C3").Value = $c2
$summary.Range("D3").Value = $d2

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 1.91

$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 2. Insert a brand new worksheet named "2022-Q3" right before "2022-Q2"
#    (i.e. right after "总计") and fill in the fund holdings table.
# ---------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($existingQ2)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "001481"
$q3.Range("C2").Value = "华宝油气（QDII）美元"
$q3.Range("D2").Value = "45.98"
$q3.Range("E2").Value = "94.53"
$q3.Range("F2").Value = "2.08"
$q3.Range("G2").Value = "0.9564"
$q3.Range("H2").Value = 7

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "162411"
$q3.Range("C3").Value = "华宝油气（QDII）人民币A"
$q3.Range("D3").Value = "28.25"
$q3.Range("E3").Value = "94.53"
$q3.Range("F3").Value = "2.08"
$q3.Range("G3").Value = "0.5876"
$q3.Range("H3").Value = 7

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "007844"
$q3.Range("C4").Value = "华宝油气（QDII）人民币 C"
$q3.Range("D4").Value = "17.73"
$q3.Range("E4").Value = "94.53"
$q3.Range("F4").Value = "2.08"
$q3.Range("G4").Value = "0.3688"
$q3.Range("H4").Value = 7

# Copy header/index-column formatting from the neighbouring "2022-Q2"
# sheet so the new sheet matches its look (bold/centered headers + A
# column style).
$existingQ2.Range("B1:H1").Copy()
$q3.Range("B1").PasteSpecial(-4122)
$existingQ2.Range("A2:A4").Copy()
$q3.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Restore "2021-Q2" as the active/selected sheet (it was the active
#    sheet before this edit).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
